# Applies the edit described by the diff to sheet1 ("Dados Pubmatic"):
#  - B2 changes from 11 to -1.7
#  - Columns D, E, G (rows 2-13) are shifted down by one row (a blank cell
#    is effectively inserted at D2/E2/G2, pushing the old row-13 values
#    into a brand-new row 14), and D3 additionally gets an explicit
#    override of -6.5 (instead of the plain shifted-down value).
#  - The new row 14 ends up with present-but-empty D14/E14/G14 cells.
#  - The sheet selection moves from D18 to H20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B2 explicit value change -------------------------------------------
$ws.Range("B2").Value = -1.7

# --- Row 2: D2, E2, G2 become empty (shifted away) ----------------------
$ws.Range("D2").Clear()
$ws.Range("E2").Clear()
$ws.Range("G2").Clear()

# --- Columns D, E, G for rows 3-13: new value = old value from the row
#     above (i.e. shift down by one), except D3 which is further
#     overridden to -6.5. -----------------------------------------------
$dValues = @{3="-6.5"; 4="-4.1"; 5="-7.8"; 6="5"; 7="11.9"; 8="10.6"; 9="3.6"; 10="4"; 11="-7"; 12="1"; 13="-8.8"}
$eValues = @{3="51"; 4="59.5"; 5="59.5"; 6="78"; 7="62"; 8="70"; 9="66"; 10="88"; 11="62"; 12="68"; 13="63.5"}
$gValues = @{3="-7.1"; 4="4.7"; 5="-10.3"; 6="5"; 7="6.6"; 8="1.8"; 9="-4.1"; 10="-2.2"; 11="-6.3"; 12="1.1"; 13="-10.1"}

foreach ($r in 3..13) {
    $ws.Range("D$r").Value = [double]$dValues[[string]$r]
    $ws.Range("E$r").Value = [double]$eValues[[string]$r]
    $ws.Range("G$r").Value = [double]$gValues[[string]$r]
}

# --- New row 14: D14/E14/G14 exist but stay empty ------------------------
# Touching the font (even to its own default) is enough to make the
# otherwise-empty cell persist in the saved sheet.
$ws.Range("D14").Font.Name = "Calibri"
$ws.Range("E14").Font.Name = "Calibri"
$ws.Range("G14").Font.Name = "Calibri"

# --- Selection moves to H20 ----------------------------------------------
$ws.Range("H20").Select()
